$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = 0.2616422292922379
$ws.Range("C2").Value = 1.03031647517663
$ws.Range("D2").Value = 3.375314179760817
$ws.Range("E2").Value = 1.837202813997632
$ws.Range("F2").Value = 1.863384125668241

# Row 3 (Q0)
$ws.Range("B3").Value = -0.689055782712207
$ws.Range("C3").Value = 1.586502696976844
$ws.Range("D3").Value = 13.20147267346463
$ws.Range("E3").Value = 3.633383089279828
$ws.Range("F3").Value = 3.651397942595873
$ws.Range("G3").Value = 22

# Row 4 (Q1)
$ws.Range("B4").Value = 0.4613554873173039
$ws.Range("C4").Value = 1.13843404674352
$ws.Range("D4").Value = 5.247316402875179
$ws.Range("E4").Value = 2.290702163720805
$ws.Range("F4").Value = 2.299171784155603
$ws.Range("G4").Value = 21
